$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the timestamps in column A (rows 2 to 97) forward by 5 days
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 5
}

# Update the Actual Production values in column B (rows 27 to 40)
$newValues = @{
    27 = 18
    28 = 73
    29 = 159
    30 = 273
    31 = 401
    32 = 538
    33 = 680
    34 = 777
    35 = 921
    36 = 1069
    37 = 1209
    38 = 1346
    39 = 1421
    40 = 1529
}

foreach ($r in $newValues.Keys) {
    $ws.Cells.Item($r, 2).Value = $newValues[$r]
}
